$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.368.98"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "3.310.12"
$ws.Range("E3").Value = "  -3.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'577.34"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").Value = "'172.74"
$ws.Range("E6").Value = "  -8.94%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "3.307.66"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("E10").Value = "  -7.01%  "
$ws.Range("D11").Value = "'0.572"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("D12").Value = "'45.05"
$ws.Range("E12").Value = "  -6.35%  "
$ws.Range("E13").Value = "  -5.98%  "
$ws.Range("D14").Value = "'667.62"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").Value = "3.848.58"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").Value = "'8.33"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "67.490.82"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.118"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.315.45"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "'17.32"
$ws.Range("E20").Value = "  -5.05%  "
$ws.Range("D21").Value = "'10.85"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'0.882"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("D23").Value = "'5.38"
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("D24").Value = "'16.87"
$ws.Range("E24").Value = "  -6.70%  "
$ws.Range("D25").Value = "'97.60"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").Value = "'3.81"
$ws.Range("E26").Value = "  -6.93%  "
$ws.Range("E27").Value = "  -7.83%  "
$ws.Range("E28").Value = "  -6.63%  "
$ws.Range("D29").Value = "'33.30"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "'8.35"
$ws.Range("E30").Value = "  -5.48%  "
$ws.Range("D31").Value = "'7.27"
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("D32").Value = "'586.75"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "'10.88"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Value = "3.694.52"
$ws.Range("E36").Value = "  -8.92%  "
$ws.Range("D37").Value = "'56.55"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -18.40%  "
$ws.Range("D39").Value = "'0.130"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").Value = "'32.58"
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("E41").Value = "  -8.41%  "
$ws.Range("D42").Value = "'3.07"
$ws.Range("E42").Value = "  -7.75%  "
$ws.Range("E43").Value = "  -5.43%  "
$ws.Range("D44").Value = "'3.27"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").Value = "0.0₃0653"
$ws.Range("E45").Value = "  -8.98%  "
$ws.Range("D46").Value = "'0.0403"
$ws.Range("E46").Value = "  -6.24%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.57"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.127"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -4.24%  "
$ws.Range("D51").Value = "'126.77"
$ws.Range("E51").Value = "  -2.24%  "
